$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1618.2307
$ws.Range("I19").Value = 786
$ws.Range("J19").Value = 2138.375
$ws.Range("K19").Value = 786
$ws.Range("L19").Value = 2138.375
$ws.Range("M19").Value = -611
$ws.Range("N19").Value = -2488.375

$ws.Range("H41").Value = 3214.158
$ws.Range("I41").Value = 3243.8333
$ws.Range("K41").Value = 3243.8333
$ws.Range("M41").Value = -2803.8333

$ws.Range("H69").Value = 43126
$ws.Range("J69").Value = 54002
$ws.Range("L69").Value = 162006
$ws.Range("N69").Value = -163754

$ws.Range("H70").Value = 2399.8
$ws.Range("I70").Value = 2166.4443
$ws.Range("K70").Value = 6499.3329
$ws.Range("M70").Value = -6229.3329

$ws.Range("H72").Value = 43126
$ws.Range("J72").Value = 54002
$ws.Range("L72").Value = 486018
$ws.Range("N72").Value = -494754

$ws.Range("H73").Value = 2399.8
$ws.Range("I73").Value = 2166.4443
$ws.Range("K73").Value = 6499.3329
$ws.Range("M73").Value = -5563.3329

$ws.Range("H74").Value = 6020
$ws.Range("I74").Value = 6020
$ws.Range("K74").Value = 6020
$ws.Range("M74").Value = -5084

$ws.Range("H76").Value = 3761.3635
$ws.Range("J76").Value = 3380
$ws.Range("L76").Value = 3380
$ws.Range("N76").Value = -4010

$ws.Range("H77").Value = 6020
$ws.Range("I77").Value = 6020
$ws.Range("K77").Value = 30100
$ws.Range("M77").Value = -25420

$ws.Range("H79").Value = 3761.3635
$ws.Range("J79").Value = 3380
$ws.Range("L79").Value = 3380
$ws.Range("N79").Value = -5564

$ws.Range("H137").Value = 4149
$ws.Range("I137").Value = 2523.9
$ws.Range("K137").Value = 7571.700000000001
$ws.Range("M137").Value = -5021.700000000001

$ws.Range("H138").Value = 3339.0122
$ws.Range("J138").Value = 3424.3816
$ws.Range("L138").Value = 10273.1448
$ws.Range("N138").Value = -20553.1448

$ws.Range("H141").Value = 3503.4167
$ws.Range("I141").Value = 3513
$ws.Range("K141").Value = 10539
$ws.Range("M141").Value = -5359

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9436105
$ws.Range("I32").Value = 10205338
$ws.Range("K32").Value = 10205338
$ws.Range("M32").Value = -10205051

$ws.Range("H88").Value = 2132.8286
$ws.Range("I88").Value = 1293.2
$ws.Range("J88").Value = 2468.68
$ws.Range("K88").Value = 1293.2
$ws.Range("L88").Value = 2468.68
$ws.Range("M88").Value = -887.2
$ws.Range("N88").Value = -3280.68

$ws.Range("H91").Value = 2132.8286
$ws.Range("I91").Value = 1293.2
$ws.Range("J91").Value = 2468.68
$ws.Range("K91").Value = 1293.2
$ws.Range("L91").Value = 2468.68
$ws.Range("M91").Value = 110.8
$ws.Range("N91").Value = -5276.68

$ws.Range("H102").Value = 7464.467
$ws.Range("I102").Value = 7585.25
$ws.Range("K102").Value = 7585.25
$ws.Range("M102").Value = -5963.25

$ws.Range("H122").Value = 2758.7585
$ws.Range("I122").Value = 2392.6667
$ws.Range("J122").Value = 4516
$ws.Range("K122").Value = 7178.000100000001
$ws.Range("L122").Value = 13548
$ws.Range("M122").Value = -4728.000100000001
$ws.Range("N122").Value = -18448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2549.4849
$ws.Range("I20").Value = 2859.6667
$ws.Range("J20").Value = 1153.6666
$ws.Range("K20").Value = 2859.6667
$ws.Range("L20").Value = 1153.6666
$ws.Range("M20").Value = -2612.6667
$ws.Range("N20").Value = -1647.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 721.9375
$ws.Range("I22").Value = 670.06665
$ws.Range("K22").Value = 670.06665
$ws.Range("M22").Value = -320.06665

$ws.Range("H107").Value = 1571.65
$ws.Range("I107").Value = 1288
$ws.Range("J107").Value = 1803.7273
$ws.Range("K107").Value = 1288
$ws.Range("L107").Value = 1803.7273
$ws.Range("M107").Value = 632
$ws.Range("N107").Value = -5643.7273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 116.13158
$ws.Range("I2").Value = 83.882355
$ws.Range("J2").Value = 142.2381
$ws.Range("K2").Value = 503.29413
$ws.Range("L2").Value = 853.4286
$ws.Range("M2").Value = -390.29413
$ws.Range("N2").Value = -1079.4286

$ws.Range("H117").Value = 2333.1667
$ws.Range("J117").Value = 2666.3333
$ws.Range("L117").Value = 7998.999899999999
$ws.Range("N117").Value = -14882.9999

$ws.Range("H132").Value = 2207.2856
$ws.Range("J132").Value = 2335.7693
$ws.Range("L132").Value = 21021.9237
$ws.Range("N132").Value = -26081.9237

$ws.Range("H141").Value = 6205.65
$ws.Range("I141").Value = 3274.2
$ws.Range("K141").Value = 9822.599999999999
$ws.Range("M141").Value = -4642.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 21366.4
$ws.Range("J96").Value = 21366.4
$ws.Range("L96").Value = 21366.4
$ws.Range("N96").Value = -26858.4

$ws.Range("H102").Value = 9348.583000000001
$ws.Range("I102").Value = 7465.3335
$ws.Range("K102").Value = 7465.3335
$ws.Range("M102").Value = -5843.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H40").Value = 3675.0588
$ws.Range("I40").Value = 2706.3333
$ws.Range("K40").Value = 2706.3333
$ws.Range("M40").Value = -2570.3333

$ws.Range("H43").Value = 1717408.5
$ws.Range("J43").Value = 1821875
$ws.Range("L43").Value = 1821875
$ws.Range("N43").Value = -1822261

$ws.Range("H68").Value = 1000
$ws.Range("J68").Value = 1000
$ws.Range("L68").Value = 1000
$ws.Range("N68").Value = -2498

$ws.Range("H71").Value = 1000
$ws.Range("J71").Value = 1000
$ws.Range("L71").Value = 5000
$ws.Range("N71").Value = -12488

$ws.Range("H93").Value = 90921970
$ws.Range("I93").Value = 90921970
$ws.Range("K93").Value = 90921970
$ws.Range("M93").Value = -90920722

$ws.Range("H122").Value = 6279.75
$ws.Range("I122").Value = 6282.591
$ws.Range("K122").Value = 18847.773
$ws.Range("M122").Value = -16397.773

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 99995.39999999999
$ws.Range("J75").Value = 99995.39999999999
$ws.Range("L75").Value = 99995.39999999999
$ws.Range("N75").Value = -101867.4

$ws.Range("H78").Value = 99995.39999999999
$ws.Range("J78").Value = 99995.39999999999
$ws.Range("L78").Value = 299986.2
$ws.Range("N78").Value = -309346.2

$ws.Range("H113").Value = 521.6667
$ws.Range("I113").Value = 944
$ws.Range("K113").Value = 2832
$ws.Range("M113").Value = -662

$ws.Range("H122").Value = 5290.3
$ws.Range("I122").Value = 5414
$ws.Range("J122").Value = 5001.6665
$ws.Range("K122").Value = 16242
$ws.Range("L122").Value = 15004.9995
$ws.Range("M122").Value = -13792
$ws.Range("N122").Value = -19904.9995
